{"js": "// Office.js (Word JavaScript API) implementation of the diff:\n//  1. Insert a new \"Meta description\" paragraph right after the Heading1\n//     title paragraph (bold \"Meta description\" run + plain-text run with\n//     the description).\n//  2. Remove the paragraph near the end of the document that duplicated the\n//     bold title text.\n//  3. Replace the text of the trailing italic paragraph (old \"Explore the\n//     Derby Dash...\" blurb) with the new image-prompt text, keeping the\n//     italic run formatting intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst titleText = \"Play Derby Dash Free: Exciting Gameplay & Impressive RTP\";\nconst oldPromptText =\n  \"Explore the Derby Dash online slot game with nudges, free spins, multipliers, and wilds. Enjoy impeccable graphics and an RTP of 96%. Play free now!\";\nconst newPromptText =\n  \"Create a feature image for Derby Dash that features a happy Maya warrior with glasses in a cartoon style. The image should show the warrior excitedly cheering on a racehorse, with the Derby Dash logo and win symbols in the background. The image should be bright and colorful, showcasing the excitement and energy of horseracing and online slot games. Make sure to include elements that represent the game's features, such as free spins and multiplier symbols.\";\n\n// ---- Step 1: insert the \"Meta description\" paragraph right after the\n// Heading1 title (i.e. right before the first paragraph that follows it). ----\nconst headingParagraph = paragraphs.items[0];\nconst afterHeading = paragraphs.items[1];\nconst insertionPoint = afterHeading.getRange(\"Start\");\n\nconst metaDescriptionSuffix =\n  \": Explore the Derby Dash online slot game with nudges, free spins, multipliers, and wilds. Enjoy impeccable graphics and an RTP of 96%. Play free now!\";\n\n// insertOoxml requires a Flat-OPC wrapped package. We provide the target\n// paragraph plus a trailing empty paragraph so Word performs an actual\n// paragraph-break insert (a single <w:p> would merge into the neighbouring\n// paragraph instead of creating a new one); the extra empty paragraph is\n// deleted again right afterwards.\nconst flatOpc =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  \"<pkg:xmlData>\" +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  \"<w:body>\" +\n  \"<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>\" +\n  metaDescriptionSuffix +\n  \"</w:t></w:r></w:p>\" +\n  \"<w:p/>\" +\n  \"</w:body>\" +\n  \"</w:document>\" +\n  \"</pkg:xmlData>\" +\n  \"</pkg:part>\" +\n  \"</pkg:package>\";\n\ninsertionPoint.insertOoxml(flatOpc, \"Before\");\nawait context.sync();\n\n// Clean up the stray empty paragraph that insertOoxml leaves behind right\n// after the inserted \"Meta description\" paragraph.\nconst afterInsertParagraphs = body.paragraphs;\nafterInsertParagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < afterInsertParagraphs.items.length; i++) {\n  const p = afterInsertParagraphs.items[i];\n  if (i > 0 && p.text === \"\") {\n    p.delete();\n    break;\n  }\n}\nawait context.sync();\n\n// ---- Step 2 & 3: near the end of the document, delete the duplicated bold\n// title paragraph and rewrite the italic paragraph's text. ----\nconst finalParagraphs = body.paragraphs;\nfinalParagraphs.load(\"text\");\nawait context.sync();\n\nlet duplicateTitleIndex = -1;\nlet italicParagraphIndex = -1;\nconst items = finalParagraphs.items;\nfor (let i = items.length - 1; i >= 1; i--) {\n  const t = items[i].text;\n  if (italicParagraphIndex === -1 && t === oldPromptText) {\n    italicParagraphIndex = i;\n    continue;\n  }\n  if (duplicateTitleIndex === -1 && t === titleText) {\n    duplicateTitleIndex = i;\n  }\n  if (duplicateTitleIndex !== -1 && italicParagraphIndex !== -1) {\n    break;\n  }\n}\n\nif (duplicateTitleIndex !== -1) {\n  items[duplicateTitleIndex].delete();\n  await context.sync();\n}\n\nif (italicParagraphIndex !== -1) {\n  const italicRange = items[italicParagraphIndex].getRange();\n  italicRange.insertText(newPromptText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop implementation of the diff:\n#  1. Insert a new \"Meta description\" paragraph right after the Heading1\n#     title paragraph (bold \"Meta description\" run + plain-text run with\n#     the description).\n#  2. Remove the paragraph near the end of the document that duplicated the\n#     bold title text.\n#  3. Replace the text of the trailing italic paragraph (old \"Explore the\n#     Derby Dash...\" blurb) with the new image-prompt text, keeping the\n#     italic run formatting intact.\n\n$d = $word.ActiveDocument\n\n$titleText = \"Play Derby Dash Free: Exciting Gameplay & Impressive RTP\"\n$oldPromptText = \"Explore the Derby Dash online slot game with nudges, free spins, multipliers, and wilds. Enjoy impeccable graphics and an RTP of 96%. Play free now!\"\n$newPromptText = \"Create a feature image for Derby Dash that features a happy Maya warrior with glasses in a cartoon style. The image should show the warrior excitedly cheering on a racehorse, with the Derby Dash logo and win symbols in the background. The image should be bright and colorful, showcasing the excitement and energy of horseracing and online slot games. Make sure to include elements that represent the game's features, such as free spins and multiplier symbols.\"\n\n# ---- Step 1: insert the \"Meta description\" paragraph right after the\n# Heading1 title paragraph. ----\n$titlePara = $d.Paragraphs.Item(1)\n$titleRange = $titlePara.Range\n$titleRange.Collapse(0)              # wdCollapseEnd\n$titleRange.InsertParagraphAfter()   # creates a new (initially empty) paragraph right after the title\n\n# Fill that freshly-created empty paragraph via an OOXML (WordprocessingML)\n# fragment so we get an exact bold \"Meta description\" run followed by the\n# plain description run, matching the target markup precisely.\n$metaSuffix = \": Explore the Derby Dash online slot game with nudges, free spins, multipliers, and wilds. Enjoy impeccable graphics and an RTP of 96%. Play free now!\"\n\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' +\n    '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>' + $metaSuffix + '</w:t></w:r></w:p>' +\n    '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n\n$newMetaPara = $d.Paragraphs.Item(2)\n$newMetaPara.Range.InsertXML($flatOpc)\n\n# ---- Step 2 & 3: near the end of the document, delete the duplicated bold\n# title paragraph and rewrite the italic paragraph's text. ----\n$count = $d.Paragraphs.Count\n\n$duplicateTitleIndex = -1\n$italicParagraphIndex = -1\n\nfor ($i = $count; $i -ge 2; $i--) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($italicParagraphIndex -eq -1 -and $t -eq $oldPromptText) {\n        $italicParagraphIndex = $i\n        continue\n    }\n    if ($duplicateTitleIndex -eq -1 -and $t -eq $titleText) {\n        $duplicateTitleIndex = $i\n    }\n    if ($duplicateTitleIndex -ne -1 -and $italicParagraphIndex -ne -1) {\n        break\n    }\n}\n\nif ($duplicateTitleIndex -ne -1) {\n    $d.Paragraphs.Item($duplicateTitleIndex).Range.Delete()\n    # deleting a paragraph before the italic one shifts its index down by one\n    if ($italicParagraphIndex -gt $duplicateTitleIndex) {\n        $italicParagraphIndex = $italicParagraphIndex - 1\n    }\n}\n\nif ($italicParagraphIndex -ne -1) {\n    $italicPara = $d.Paragraphs.Item($italicParagraphIndex)\n    $s = $italicPara.Range.Start\n    $e = $italicPara.Range.End\n    $replaceRange = $d.Range($s, $e)\n    $replaceRange.Text = $newPromptText\n}\n"}
